$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date for all existing data rows (2-155)
# from 45190 to 45192.
$lastRow = 155
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}

# Row 155 ends up with an explicit row height in the target file.
$ws.Rows.Item($lastRow).RowHeight = 15

# Append a new row (156) with the new record.
$newRow = 156

$ws.Cells.Item($newRow, 1).Value2 = "A 44619-2023"

$ws.Cells.Item($newRow, 2).Value2 = 45189
$ws.Cells.Item($newRow, 2).NumberFormat = $ws.Cells.Item($lastRow, 2).NumberFormat

$ws.Cells.Item($newRow, 3).Value2 = 45192
$ws.Cells.Item($newRow, 3).NumberFormat = $ws.Cells.Item($lastRow, 3).NumberFormat

$ws.Cells.Item($newRow, 4).Value2 = "STOCKHOLMS LÄN"
$ws.Cells.Item($newRow, 5).Value2 = "NYNÄSHAMN"

$ws.Cells.Item($newRow, 7).Value2 = 7.6
$ws.Cells.Item($newRow, 8).Value2 = 0
$ws.Cells.Item($newRow, 9).Value2 = 0
$ws.Cells.Item($newRow, 10).Value2 = 0
$ws.Cells.Item($newRow, 11).Value2 = 0
$ws.Cells.Item($newRow, 12).Value2 = 0
$ws.Cells.Item($newRow, 13).Value2 = 0
$ws.Cells.Item($newRow, 14).Value2 = 0
$ws.Cells.Item($newRow, 15).Value2 = 0
$ws.Cells.Item($newRow, 16).Value2 = 0
$ws.Cells.Item($newRow, 17).Value2 = 0

$ws.Cells.Item($newRow, 18).WrapText = $ws.Cells.Item($lastRow, 18).WrapText

$ws.Rows.Item($newRow).RowHeight = 15
